$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.17831799999999
$ws.Range("H2").Value = 288.534954
$ws.Range("I2").Value = 0.7237598617297997
$ws.Range("J2").Value = 0.7237598617297996
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 223.5919553858307
$ws.Range("R2").Value = 2012.327598472476
$ws.Range("S2").Value = 0.0323804661464409
$ws.Range("T2").Value = 0.03238046614644088

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.17831799999999
$ws.Range("H3").Value = 288.534954
$ws.Range("I3").Value = 0.7237598617297997
$ws.Range("J3").Value = 0.7237598617297996
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 1853.82935312512
$ws.Range("R3").Value = 16684.46417812608
$ws.Range("S3").Value = 0.2684705650816561
$ws.Range("T3").Value = 0.268470565081656

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.17831799999999
$ws.Range("H4").Value = 288.534954
$ws.Range("I4").Value = 0.7237598617297997
$ws.Range("J4").Value = 0.7237598617297996
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 2920.248644172283
$ws.Range("R4").Value = 26282.23779755055
$ws.Range("S4").Value = 0.4229088305017027
$ws.Range("T4").Value = 0.4229088305017027

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 13.23504133333333
$ws.Range("H5").Value = 39.705124
$ws.Range("I5").Value = 0.09959616558694152
$ws.Range("J5").Value = 0.0995961655869415
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.324764666666667
$ws.Range("N5").Value = 6.974294
$ws.Range("O5").Value = 0.04473923998638302
$ws.Range("P5").Value = 0.04473923998638301
$ws.Range("Q5").Value = 30.76835645360622
$ws.Range("R5").Value = 276.915208082456
$ws.Range("S5").Value = 0.004455856753917718
$ws.Range("T5").Value = 0.004455856753917717

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.23504133333333
$ws.Range("H6").Value = 39.705124
$ws.Range("I6").Value = 0.09959616558694152
$ws.Range("J6").Value = 0.0995961655869415
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.27491966666667
$ws.Range("N6").Value = 57.824759
$ws.Range("O6").Value = 0.3709387315842666
$ws.Range("P6").Value = 0.3709387315842665
$ws.Range("Q6").Value = 255.1043584850129
$ws.Range("R6").Value = 2295.939226365116
$ws.Range("S6").Value = 0.03694407533347667
$ws.Range("T6").Value = 0.03694407533347666

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.23504133333333
$ws.Range("H7").Value = 39.705124
$ws.Range("I7").Value = 0.09959616558694152
$ws.Range("J7").Value = 0.0995961655869415
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.5843220284293504
$ws.Range("P7").Value = 0.5843220284293504
$ws.Range("Q7").Value = 401.8536850398111
$ws.Range("R7").Value = 3616.6831653583
$ws.Range("S7").Value = 0.05819623349954714
$ws.Range("T7").Value = 0.05819623349954712

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.47369766666667
$ws.Range("H8").Value = 70.421093
$ws.Range("I8").Value = 0.1766439726832589
$ws.Range("J8").Value = 0.1766439726832589
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.324764666666667
$ws.Range("N8").Value = 6.974294
$ws.Range("O8").Value = 0.04473923998638302
$ws.Range("P8").Value = 0.04473923998638301
$ws.Range("Q8").Value = 54.57082293148245
$ws.Range("R8").Value = 491.137406383342
$ws.Range("S8").Value = 0.007902917086024407
$ws.Range("T8").Value = 0.007902917086024405

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.47369766666667
$ws.Range("H9").Value = 70.421093
$ws.Range("I9").Value = 0.1766439726832589
$ws.Range("J9").Value = 0.1766439726832589
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.27491966666667
$ws.Range("N9").Value = 57.824759
$ws.Range("O9").Value = 0.3709387315842666
$ws.Range("P9").Value = 0.3709387315842665
$ws.Range("Q9").Value = 452.4536368046208
$ws.Range("R9").Value = 4072.082731241587
$ws.Range("S9").Value = 0.06552409116913391
$ws.Range("T9").Value = 0.06552409116913388

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.47369766666667
$ws.Range("H10").Value = 70.421093
$ws.Range("I10").Value = 0.1766439726832589
$ws.Range("J10").Value = 0.1766439726832589
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.5843220284293504
$ws.Range("P10").Value = 0.5843220284293504
$ws.Range("Q10").Value = 712.7285568124972
$ws.Range("R10").Value = 6414.557011312475
$ws.Range("S10").Value = 0.1032169644281006
$ws.Range("T10").Value = 0.1032169644281006

